# "cambios de las fracciones"
# Update the reporting-period / validation dates in row 8 and move the
# active cell selection from C18 to C14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Fecha de inicio del periodo que se informa (B8): 1/1/2022 -> 7/1/2022
$ws.Range("B8").Value = "7/1/2022"

# Fecha de término del periodo que se informa (C8): 6/30/2022 -> 12/31/2022
$ws.Range("C8").Value = "12/31/2022"

# Fecha de validación (I8): 7/11/2022 -> 1/10/2023
$ws.Range("I8").Value = "1/10/2023"

# Fecha de actualización (J8): 7/11/2022 -> 1/10/2023
$ws.Range("J8").Value = "1/10/2023"

# Move the saved selection from C18 to C14
$ws.Range("C14").Select()
